$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.329.36'
$ws.Range("E2").Value = '  +4.99%  '
$ws.Range("D3").Value = '2.466.58'
$ws.Range("E3").Value = '  +6.22%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.72%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +2.55%  '
$ws.Range("D9").Value = '2.465.37'
$ws.Range("E9").Value = '  +6.25%  '
$ws.Range("E10").Value = '  +5.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.74'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.10%  '
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("E13").Value = '  +5.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +12.74%  '
$ws.Range("D15").Value = '2.909.93'
$ws.Range("E15").Value = '  +5.92%  '
$ws.Range("D16").Value = '63.257.74'
$ws.Range("E16").Value = '  +4.70%  '
$ws.Range("E17").Value = '  +7.47%  '
$ws.Range("D18").Value = '2.470.33'
$ws.Range("E18").Value = '  +5.96%  '
$ws.Range("E19").Value = '  +6.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.73%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.54%  '
$ws.Range("E25").Value = '  +1.94%  '
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("E27").Value = '  +8.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("E29").Value = '  +9.26%  '
$ws.Range("D30").Value = '0.0₃0818'
$ws.Range("E30").Value = '  +13.05%  '
$ws.Range("E31").Value = '  +14.41%  '
$ws.Range("E32").Value = '  +7.38%  '
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.98%  '
$ws.Range("E35").Value = '  +4.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '371.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +16.21%  '
$ws.Range("E38").Value = '  +9.01%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.73'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.43%  '
$ws.Range("E46").Value = '  +5.71%  '
$ws.Range("E47").Value = '  +2.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0521'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.98%  '
$ws.Range("D49").Value = '0.0₆0241'
$ws.Range("E49").Value = '  +10.91%  '
$ws.Range("E50").Value = '  +5.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.74%  '
